# daily auto push: 2025-10-05 07:23 UTC
# Append the newest log entry as a new row at the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (e.g. "2025/10/05"), matching the
# existing rows above. Force text so Excel doesn't auto-convert the
# date-looking string into a date serial number, then drop the temporary
# number format so no stray cell style is left behind.
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "2025/10/05"
$ws.Range("A65").ClearFormats()

$ws.Range("B65").Value = "日"
$ws.Range("C65").Value = 16
$ws.Range("D65").Value = 201
